# Apply "Updated symbol list on Sun Dec 18 22:46:03 UTC 2022 with GitHub Actions"
# - Refreshed cryptocurrency prices in column D (stored as text, matching the
#   source sheet's existing inline-string convention).
# - A few "Bestin24h"/"Worstin24h" suffix tweaks in column E.
# - Rows 42/43: CEJI and BKEXToken swapped positions (columns B, C, D, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: for cells whose new text looks like a plain number,
# prefix with an apostrophe so Excel keeps storing it as TEXT (matching the
# original cell type), then reset the cell style so no stray 'quote prefix'
# formatting is introduced.

$ws.Range("D2").Value = "'251.13"
$ws.Range("D2").Style = "Normal"

$ws.Range("D4").Value = "'5.556"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = "'0.05692"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'6.449"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'0.8090"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'1.038"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.1431"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.07283"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.03142"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'0.02919"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'0.09268"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.001658"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'3.210"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'0.04747"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.0005810"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16OneONE"

$ws.Range("D18").Value = "'0.006444"
$ws.Range("D18").Style = "Normal"

$ws.Range("D19").Value = "'0.005075"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"

$ws.Range("D20").Value = "'0.001048"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").Value = "'0.0001497"
$ws.Range("D21").Style = "Normal"

$ws.Range("D22").Value = "'3.990"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'3.371"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").Value = "'2.112"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").Value = "'0.3323"
$ws.Range("D25").Style = "Normal"

$ws.Range("D27").Value = "'0.0003094"
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").Value = "'0.04133"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006876"
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1045"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003194"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.009619"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005636"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'0.7838"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

$ws.Range("D48").Value = "'0.01683"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"

$ws.Range("D49").Value = "'0.00002096"
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").Value = "'0.01008"
$ws.Range("D50").Style = "Normal"
